$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = "Tax being charged to invoice 123456789 on ship to 987654321. BU 10/04/2017 ***DUPLICATE SR***"
$ws.Range("B3").Value = "Lorem ipsum dolor sit amet, consectetur adipiscing elit, sed do eiusmod tempor incididunt ut labore et dolore magna aliqua. Ut enim ad minim veniam, quis nostrud exercitation ullamco laboris nisi ut aliquip ex ea commodo consequat. Duis aute irure dolor in reprehenderit in voluptate velit esse cillum dolore eu fugiat nulla pariatur. Excepteur sint occaecat cupidatat non proident, sunt in culpa qui officia deserunt mollit anim id est laborum"
$ws.Range("A1").Value = "Service Request Number"
$ws.Range("B2").Value = "CSC 08/08/16 ABCD - wrong product on four invoices 12345, 67890, 09876, 54321"

$ws.Range("A2").Value = 12345
$ws.Range("A3").Value = 67890
$ws.Range("A4").Value = 54321

$ws.Range("C9").Select()
